# Bishop State Community College Organizations workbook — column rework:
#   * "Organization Name" / "Categories" swap position and the category
#     column is renamed/relocated to column A ("Category"), organization
#     name moves to column B.
#   * Several headers are renamed (Org URL -> Organization Link, Image URL ->
#     Logo Link, Phone -> Phone Number, Website dropped, LinkedIn/Instagram/
#     Facebook/Twitter -> *_Link and shifted left into the old Website slot).
#   * Two new trailing columns are appended: "Youtube Link" and "Tiktok Link".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: re-derive every data row (rows 2-32) from the OLD layout
#         (A:Name  B:Category  C:Url  D:Image  E:Desc  F:Email  G:Phone
#          H:Website  I:LinkedIn  J:Instagram  K:Facebook  L:Twitter)
#         into the NEW layout
#         (A:Category B:Name C:Url D:Logo E:Desc F:Email G:Phone
#          H:LinkedinLink I:InstagramLink J:FacebookLink K:TwitterLink
#          L:YoutubeLink(new,empty) M:TiktokLink(new,empty))
# ---------------------------------------------------------------------
for ($r = 2; $r -le 32; $r++) {
    $old = $ws.Range("A$r`:L$r").Value2

    $new = New-Object 'object[,]' 1,13
    $new[0,0]  = $old[1,2]    # Category      <- old Categories (B)
    $new[0,1]  = $old[1,1]    # Org Name      <- old Organization Name (A)
    $new[0,2]  = $old[1,3]    # Org Link      <- old Org URL (C)
    $new[0,3]  = $old[1,4]    # Logo Link     <- old Image URL (D)
    $new[0,4]  = $old[1,5]    # Description   <- old Description (E)
    $new[0,5]  = $old[1,6]    # Email         <- old Email (F)
    $new[0,6]  = $old[1,7]    # Phone Number  <- old Phone (G)
    $new[0,7]  = $old[1,9]    # Linkedin Link <- old LinkedIn (I)
    $new[0,8]  = $old[1,10]   # Instagram Link<- old Instagram (J)
    $new[0,9]  = $old[1,11]   # Facebook Link <- old Facebook (K)
    $new[0,10] = $old[1,12]   # Twitter Link  <- old Twitter (L)
    $new[0,11] = $null        # Youtube Link  (new, empty)
    $new[0,12] = $null        # Tiktok Link   (new, empty)

    $ws.Range("A$r`:M$r").Value2 = $new
}

# ---------------------------------------------------------------------
# Step 2: rewrite the header row (row 1) with the new titles/order.
# ---------------------------------------------------------------------
$ws.Cells.Item(1,1).Value2  = "Category"
$ws.Cells.Item(1,2).Value2  = "Organization Name"
$ws.Cells.Item(1,3).Value2  = "Organization Link"
$ws.Cells.Item(1,4).Value2  = "Logo Link"
$ws.Cells.Item(1,5).Value2  = "Description"
$ws.Cells.Item(1,6).Value2  = "Email"
$ws.Cells.Item(1,7).Value2  = "Phone Number"
$ws.Cells.Item(1,8).Value2  = "Linkedin Link"
$ws.Cells.Item(1,9).Value2  = "Instagram Link"
$ws.Cells.Item(1,10).Value2 = "Facebook Link"
$ws.Cells.Item(1,11).Value2 = "Twitter Link"
$ws.Cells.Item(1,12).Value2 = "Youtube Link"
$ws.Cells.Item(1,13).Value2 = "Tiktok Link"

# Column M is brand-new, so copy the bold/bordered/centered header style
# from an existing header cell onto it (format-only paste).
$ws.Cells.Item(1,1).Copy()
$ws.Cells.Item(1,13).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 3: column widths. ColumnWidth round-trips into the saved xlsx
# width with a constant +5/6 character offset, so subtract that back
# out to land exactly on the target widths.
# ---------------------------------------------------------------------
$wOffset = 0.8333333333333334
function Set-ColWidth($colIndex, $target) {
    $ws.Columns.Item($colIndex).ColumnWidth = $target - $wOffset
}

Set-ColWidth 1  14   # A Category
Set-ColWidth 2  34   # B Organization Name
Set-ColWidth 3  50   # C Organization Link
Set-ColWidth 4  11   # D Logo Link
Set-ColWidth 5  13   # E Description
Set-ColWidth 6  7    # F Email
Set-ColWidth 7  14   # G Phone Number
Set-ColWidth 8  50   # H Linkedin Link
Set-ColWidth 9  39   # I Instagram Link
Set-ColWidth 10 38   # J Facebook Link
Set-ColWidth 11 33   # K Twitter Link
Set-ColWidth 12 14   # L Youtube Link
Set-ColWidth 13 13   # M Tiktok Link
